$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 3).Value = 324535
$ws.Cells.Item(2, 4).Value = 413482458
$ws.Cells.Item(4, 3).Value = 327
$ws.Cells.Item(4, 4).Value = 467692
$ws.Cells.Item(10, 3).Value = 118013
$ws.Cells.Item(10, 4).Value = 172917118
$ws.Cells.Item(12, 3).Value = 60152
$ws.Cells.Item(12, 4).Value = 86814505
$ws.Cells.Item(16, 3).Value = 4027
$ws.Cells.Item(16, 4).Value = 5715197
$ws.Cells.Item(20, 3).Value = 6863
$ws.Cells.Item(20, 4).Value = 9580704
$ws.Cells.Item(22, 3).Value = 78268
$ws.Cells.Item(22, 4).Value = 97545640
$ws.Cells.Item(28, 3).Value = 32688
$ws.Cells.Item(28, 4).Value = 47844999
$ws.Cells.Item(30, 3).Value = 11583
$ws.Cells.Item(30, 4).Value = 16661747
$ws.Cells.Item(33, 3).Value = 1568
$ws.Cells.Item(33, 4).Value = 2203281
$ws.Cells.Item(34, 3).Value = 14
$ws.Cells.Item(34, 4).Value = 20684
$ws.Cells.Item(35, 3).Value = 1890
$ws.Cells.Item(35, 4).Value = 2670168
$ws.Cells.Item(36, 3).Value = 98146
$ws.Cells.Item(36, 4).Value = 123446064
$ws.Cells.Item(44, 3).Value = 44651
$ws.Cells.Item(44, 4).Value = 65437560
$ws.Cells.Item(46, 3).Value = 9245
$ws.Cells.Item(46, 4).Value = 13261108
$ws.Cells.Item(48, 3).Value = 1416
$ws.Cells.Item(48, 4).Value = 1966603
$ws.Cells.Item(51, 3).Value = 2393
$ws.Cells.Item(51, 4).Value = 3343402
$ws.Cells.Item(52, 3).Value = 69764
$ws.Cells.Item(52, 4).Value = 87497943
$ws.Cells.Item(59, 3).Value = 28415
$ws.Cells.Item(59, 4).Value = 41670204
$ws.Cells.Item(62, 3).Value = 11252
$ws.Cells.Item(62, 4).Value = 16271212
$ws.Cells.Item(64, 3).Value = 1366
$ws.Cells.Item(64, 4).Value = 1909237
$ws.Cells.Item(68, 3).Value = 1526
$ws.Cells.Item(68, 4).Value = 2138114
$ws.Cells.Item(70, 3).Value = 20699
$ws.Cells.Item(70, 4).Value = 27107108
$ws.Cells.Item(74, 3).Value = 7660
$ws.Cells.Item(74, 4).Value = 11216856
$ws.Cells.Item(76, 3).Value = 5177
$ws.Cells.Item(76, 4).Value = 7517372
$ws.Cells.Item(78, 3).Value = 285
$ws.Cells.Item(78, 4).Value = 400583
$ws.Cells.Item(79, 3).Value = 142338
$ws.Cells.Item(79, 4).Value = 177403150
$ws.Cells.Item(85, 3).Value = 64050
$ws.Cells.Item(85, 4).Value = 93870725
$ws.Cells.Item(88, 3).Value = 30047
$ws.Cells.Item(88, 4).Value = 43464983
$ws.Cells.Item(90, 3).Value = 2753
$ws.Cells.Item(90, 4).Value = 3963957
$ws.Cells.Item(91, 3).Value = 2902
$ws.Cells.Item(91, 4).Value = 4103847
$ws.Cells.Item(92, 3).Value = 33775
$ws.Cells.Item(92, 4).Value = 45777271
$ws.Cells.Item(96, 3).Value = 8196
$ws.Cells.Item(96, 4).Value = 12048402
$ws.Cells.Item(98, 3).Value = 7531
$ws.Cells.Item(98, 4).Value = 10926568
$ws.Cells.Item(100, 3).Value = 540
$ws.Cells.Item(100, 4).Value = 766651
$ws.Cells.Item(101, 3).Value = 503
$ws.Cells.Item(101, 4).Value = 725891
$ws.Cells.Item(102, 3).Value = 10515
$ws.Cells.Item(102, 4).Value = 16096474
$ws.Cells.Item(104, 3).Value = 2586
$ws.Cells.Item(104, 4).Value = 4228621
$ws.Cells.Item(106, 3).Value = 3504
$ws.Cells.Item(106, 4).Value = 5735573
$ws.Cells.Item(108, 3).Value = 161
$ws.Cells.Item(108, 4).Value = 265445
$ws.Cells.Item(109, 3).Value = 204
$ws.Cells.Item(109, 4).Value = 318530
$ws.Cells.Item(110, 3).Value = 143033
$ws.Cells.Item(110, 4).Value = 176863546
$ws.Cells.Item(113, 3).Value = 7
$ws.Cells.Item(113, 4).Value = 10395
$ws.Cells.Item(116, 3).Value = 53156
$ws.Cells.Item(116, 4).Value = 77912310
$ws.Cells.Item(117, 3).Value = 86
$ws.Cells.Item(117, 4).Value = 127459
$ws.Cells.Item(118, 3).Value = 27466
$ws.Cells.Item(118, 4).Value = 39794034
$ws.Cells.Item(122, 3).Value = 2320
$ws.Cells.Item(122, 4).Value = 3260509
$ws.Cells.Item(124, 3).Value = 522361
$ws.Cells.Item(124, 4).Value = 689896786
$ws.Cells.Item(129, 3).Value = 1388
$ws.Cells.Item(129, 4).Value = 2057182
$ws.Cells.Item(131, 3).Value = 210770
$ws.Cells.Item(131, 4).Value = 309844471
$ws.Cells.Item(132, 3).Value = 410
$ws.Cells.Item(132, 4).Value = 611750
$ws.Cells.Item(134, 3).Value = 186400
$ws.Cells.Item(134, 4).Value = 271064958
$ws.Cells.Item(136, 3).Value = 33
$ws.Cells.Item(136, 4).Value = 48332
$ws.Cells.Item(137, 3).Value = 2872
$ws.Cells.Item(137, 4).Value = 4032799
$ws.Cells.Item(139, 3).Value = 6538
$ws.Cells.Item(139, 4).Value = 9236919
$ws.Cells.Item(142, 3).Value = 45326
$ws.Cells.Item(142, 4).Value = 60507242
$ws.Cells.Item(148, 3).Value = 14229
$ws.Cells.Item(148, 4).Value = 20861738
$ws.Cells.Item(149, 3).Value = 3808
$ws.Cells.Item(149, 4).Value = 5491240
$ws.Cells.Item(154, 3).Value = 396
$ws.Cells.Item(154, 4).Value = 559763
$ws.Cells.Item(155, 3).Value = 17861
$ws.Cells.Item(155, 4).Value = 23606905
$ws.Cells.Item(159, 3).Value = 7285
$ws.Cells.Item(159, 4).Value = 10599404
$ws.Cells.Item(161, 3).Value = 5082
$ws.Cells.Item(161, 4).Value = 7315496
$ws.Cells.Item(164, 3).Value = 272
$ws.Cells.Item(164, 4).Value = 389364
$ws.Cells.Item(166, 3).Value = 18714
$ws.Cells.Item(166, 4).Value = 30538370
$ws.Cells.Item(167, 3).Value = 2044
$ws.Cells.Item(167, 4).Value = 3361072
$ws.Cells.Item(168, 3).Value = 278
$ws.Cells.Item(168, 4).Value = 454589
$ws.Cells.Item(171, 3).Value = 106
$ws.Cells.Item(171, 4).Value = 181449
$ws.Cells.Item(172, 3).Value = 88709
$ws.Cells.Item(172, 4).Value = 110866653
$ws.Cells.Item(177, 3).Value = 645
$ws.Cells.Item(177, 4).Value = 950588
$ws.Cells.Item(179, 3).Value = 34138
$ws.Cells.Item(179, 4).Value = 50059477
$ws.Cells.Item(181, 3).Value = 13180
$ws.Cells.Item(181, 4).Value = 19042845
$ws.Cells.Item(183, 3).Value = 1251
$ws.Cells.Item(183, 4).Value = 1751429
$ws.Cells.Item(185, 3).Value = 1700
$ws.Cells.Item(185, 4).Value = 2389005
$ws.Cells.Item(187, 3).Value = 240882
$ws.Cells.Item(187, 4).Value = 299325312
$ws.Cells.Item(189, 3).Value = 170
$ws.Cells.Item(189, 4).Value = 245236
$ws.Cells.Item(193, 3).Value = 883
$ws.Cells.Item(193, 4).Value = 1298845
$ws.Cells.Item(195, 3).Value = 87253
$ws.Cells.Item(195, 4).Value = 127897857
$ws.Cells.Item(198, 3).Value = 33367
$ws.Cells.Item(198, 4).Value = 48032595
$ws.Cells.Item(201, 3).Value = 5156
$ws.Cells.Item(201, 4).Value = 7339796
$ws.Cells.Item(204, 3).Value = 5017
$ws.Cells.Item(204, 4).Value = 6947378
$ws.Cells.Item(207, 3).Value = 266897
$ws.Cells.Item(207, 4).Value = 330258895
$ws.Cells.Item(209, 3).Value = 257
$ws.Cells.Item(209, 4).Value = 367039
$ws.Cells.Item(214, 3).Value = 621
$ws.Cells.Item(214, 4).Value = 904378
$ws.Cells.Item(216, 3).Value = 95937
$ws.Cells.Item(216, 4).Value = 140349246
$ws.Cells.Item(219, 3).Value = 52066
$ws.Cells.Item(219, 4).Value = 75250062
$ws.Cells.Item(222, 3).Value = 4709
$ws.Cells.Item(222, 4).Value = 6609823
$ws.Cells.Item(225, 3).Value = 5950
$ws.Cells.Item(225, 4).Value = 8236605
$ws.Cells.Item(228, 3).Value = 107504
$ws.Cells.Item(228, 4).Value = 134392646
$ws.Cells.Item(235, 3).Value = 49912
$ws.Cells.Item(235, 4).Value = 73117543
$ws.Cells.Item(237, 3).Value = 12608
$ws.Cells.Item(237, 4).Value = 18127904
$ws.Cells.Item(239, 3).Value = 1902
$ws.Cells.Item(239, 4).Value = 2726882
$ws.Cells.Item(241, 3).Value = 2586
$ws.Cells.Item(241, 4).Value = 3619686
$ws.Cells.Item(242, 3).Value = 260873
$ws.Cells.Item(242, 4).Value = 329377272
$ws.Cells.Item(246, 3).Value = 16
$ws.Cells.Item(246, 4).Value = 24000
$ws.Cells.Item(248, 3).Value = 837
$ws.Cells.Item(248, 4).Value = 1228904
$ws.Cells.Item(250, 3).Value = 96661
$ws.Cells.Item(250, 4).Value = 141635122
$ws.Cells.Item(253, 3).Value = 65973
$ws.Cells.Item(253, 4).Value = 95608796
$ws.Cells.Item(255, 3).Value = 2431
$ws.Cells.Item(255, 4).Value = 3427773
$ws.Cells.Item(258, 3).Value = 4748
$ws.Cells.Item(258, 4).Value = 6673279
